# Updates the Oni pricing report:
#  - "Generated:" timestamp and "RQST By:" requester lines were regenerated
#    (new run at 11:15 AM, requested by ROMERO ONIRIA instead of ARAMIS).
#    These two lines are repeated at the top of every "page" of the report.
#  - A handful of dollar figures were refreshed to match the new pricing run
#    (each down by 336, from a material-cost line that feeds the various
#    subtotals / grand totals further down the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$generatedText = "Generated: 05/24/2024, 11:15 AM"
$rqstText = "RQST By: ROMERO ONIRIA"

# Row pairs where the "Generated:" / "RQST By:" banner is repeated throughout
# the report (one pair per report "page").
$generatedRows = @(3, 53, 99, 134, 157, 182, 191, 230)
$rqstRows = @(4, 54, 100, 135, 158, 183, 192, 231)

foreach ($r in $generatedRows) {
    $ws.Range("A$r").Value = $generatedText
}

foreach ($r in $rqstRows) {
    $ws.Range("A$r").Value = $rqstText
    # The new requester name wraps onto more lines than "RQST By: ARAMIS" did,
    # so the row grows from 42.75 to 71.25 points.
    $ws.Rows.Item($r).RowHeight = 71.25
}

# Refreshed pricing figures (all down by 336 versus the previous run).
$ws.Range("G46").Value = 27999
$ws.Range("G48").Value = 94414.54
$ws.Range("B142").Value = 154440.21
$ws.Range("B144").Value = 154440.21
$ws.Range("C151").Value = 154440.21
$ws.Range("C152").Value = 154440.21
$ws.Range("C220").Value = 34611
$ws.Range("B221").Value = 27999
$ws.Range("C225").Value = 34611
$ws.Range("B228").Value = 34611
